$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.238.49"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "2.268.28"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "495.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0960"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.11%  "
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("E11").Value = "  +3.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").Value = "2.673.66"
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.25%  "
$ws.Range("D15").Value = "54.171.33"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").Value = "2.266.81"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.36%  "
$ws.Range("E19").Value = "  +3.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "303.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.45%  "
$ws.Range("E21").Value = "  +4.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "2.378.23"
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "170.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.25%  "
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0679"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.898"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.77%  "
$ws.Range("E37").Value = "  +1.59%  "
$ws.Range("E38").Value = "  +3.52%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.372"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.70%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "126.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0897"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.86%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0487"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.545"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.82%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "237.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.67%  "
$ws.Range("B48").Value = "Polygon"
$ws.Range("C48").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.371"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0205"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.24%  "
